# Update the "取得日時" (acquired date/time) timestamps in rows 2-12 of the
# "ランサーズ" sheet from 2025-11-13 12:39:17 to 2025-11-13 12:50:59.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-11-13 12:50:59"

for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
